$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data table (rows 2-9, columns A:T) reflecting updated TPM computation.
# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric metrics
$data = @(
    @("ECs","Adm","Calcr","ECs", 3,1,12.31940633333333,36.958219,0.3801768364207869,0.3801768364207869,1,0.3333333333333333,0.037013,0.111039,0.003877153058804169,0.003877153058804169,0.4559781866156666,4.103803679541,0.001474003784215346,0.001474003784215346),
    @("ECs","Adm","Calcr","MuSCs", 3,1,12.31940633333333,36.958219,0.3801768364207869,0.3801768364207869,3,1,9.509424666666666,28.528274,0.9961228469411958,0.9961228469411958,117.1504664648895,1054.354198184006,0.3787028326365715,0.3787028326365715),
    @("FAPs","Adm","Calcr","ECs", 3,1,18.11265066666666,54.33795199999999,0.5589563363143816,0.5589563363143816,1,0.3333333333333333,0.037013,0.111039,0.003877153058804169,0.003877153058804169,0.6704035391253331,6.033631852127999,0.002167159269079276,0.002167159269079276),
    @("FAPs","Adm","Calcr","MuSCs", 3,1,18.11265066666666,54.33795199999999,0.5589563363143816,0.5589563363143816,3,1,9.509424666666666,28.528274,0.9961228469411958,0.9961228469411958,172.2408870283164,1550.167983254848,0.5567891770453023,0.5567891770453023),
    @("MuSCs","Adm","Calcr","ECs", 3,1,1.603212,4.809636,0.04947511672073613,0.04947511672073613,1,0.3333333333333333,0.037013,0.111039,0.003877153058804169,0.003877153058804169,0.059339685756,0.534057171804,0.0001918226001284954,0.0001918226001284954),
    @("MuSCs","Adm","Calcr","MuSCs", 3,1,1.603212,4.809636,0.04947511672073613,0.04947511672073613,3,1,9.509424666666666,28.528274,0.9961228469411958,0.9961228469411958,15.245623738696,137.210613648264,0.04928329412060763,0.04928329412060763),
    @("Resolving-Mac","Adm","Calcr","ECs", 3,1,0.3691416666666666,1.107425,0.01139171054409548,0.01139171054409548,1,0.3333333333333333,0.037013,0.111039,0.003877153058804169,0.003877153058804169,0.01366304050833333,0.122967364575,0.00004416740538105149,0.00004416740538105149),
    @("Resolving-Mac","Adm","Calcr","MuSCs", 3,1,0.3691416666666666,1.107425,0.01139171054409548,0.01139171054409548,3,1,9.509424666666666,28.528274,0.9961228469411958,0.9961228469411958,3.510324870494444,31.59292383444999,0.01134754313871443,0.01134754313871443)
)

for ($r = 0; $r -lt $data.Count; $r++) {
    $row = $data[$r]
    $targetRow = $r + 2
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($targetRow, $c + 1).Value = $row[$c]
    }
}
